$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Set the new mock file values in column B (data.txt for row 8, export.csv for row 9)
$ws.Range("B8").Value = "data.txt"
$ws.Range("B9").Value = "export.csv"

# Update selection to B11 (as seen in the diff)
$ws.Range("B11").Select()
